# DYNO-23: Handling Parallel Operations - updated table structure
#
# 1. ShopOrderOperations: rename "PrecedingOperation" column header to
#    "PrecedingOperationID" and update its values to reference the actual
#    OperationID of the preceding operation (instead of its old sequence flag).
# 2. Update the remembered cell selections on the ShopOrders and
#    ShopOrderOperations sheets.

$wb = $excel.ActiveWorkbook

$wsShopOrders = $wb.Worksheets.Item("ShopOrders")
$wsShopOrderOperations = $wb.Worksheets.Item("ShopOrderOperations")

# Rename the "PrecedingOperation" header to "PrecedingOperationID"
$wsShopOrderOperations.Range("F1").Value = "PrecedingOperationID"

# Update values to hold the preceding operation's actual OperationID
$wsShopOrderOperations.Range("F3").Value = 100
$wsShopOrderOperations.Range("F5").Value = 201

# Update remembered selections.
# Select ShopOrders first so ShopOrderOperations ends up as the active tab,
# matching the original workbook state.
$wsShopOrders.Range("G20").Select()
$wsShopOrderOperations.Range("E25").Select()
